$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the custom column width formatting on B:G (removes the <cols> element).
# This also leaves behind transient blank cells in row 2 (since it already has
# data in B2) for every other touched column; those get cleaned up below.
$ws.Columns("B:G").ClearFormats()

# --- Row 2: headers for the whole block ---
$ws.Range("B2").Value = "particle_block-59"
$ws.Range("G2").Value = "particle_block-137"
$ws.Range("C2:F2").Clear()

# --- Block 1: MIXT, (rows 4-8) ---
$ws.Range("B4").Value = "M/Z/O/E"
$ws.Range("C4").Value = "MIXT,"

$ws.Range("B5").Value = [double]"0"
$ws.Range("C5").Value = [double]"0"
$ws.Range("D5").Value = [double]"99.9999"

$ws.Range("B6").Value = [double]"0.625"
$ws.Range("C6").Value = [double]"0"
$ws.Range("D6").Value = [double]"99.9999"

$ws.Range("B7").Value = [double]"100000"
$ws.Range("C7").Value = [double]"7.37042E-17"
$ws.Range("D7").Value = [double]"0.33158"

$ws.Range("B8").Value = [double]"10000000"
$ws.Range("C8").Value = [double]"6.29018E-20"
$ws.Range("D8").Value = [double]"14.2801"

# --- Block 2: O16 (rows 10-14) ---
$ws.Range("B10").Value = "M/Z/O/E"
$ws.Range("C10").Value = "O16"

$ws.Range("B11").Value = [double]"0"
$ws.Range("C11").Value = [double]"0"
$ws.Range("D11").Value = [double]"99.9999"

$ws.Range("B12").Value = [double]"0.625"
$ws.Range("C12").Value = [double]"0"
$ws.Range("D12").Value = [double]"99.9999"

$ws.Range("B13").Value = [double]"100000"
$ws.Range("C13").Value = [double]"2.43272E-19"
$ws.Range("D13").Value = [double]"0.297053"

$ws.Range("B14").Value = [double]"10000000"
$ws.Range("C14").Value = [double]"6.10277E-23"
$ws.Range("D14").Value = [double]"13.8157"

# --- Block 3: N14 (rows 16-20) ---
$ws.Range("B16").Value = "M/Z/O/E"
$ws.Range("C16").Value = "N14"

$ws.Range("B17").Value = [double]"0"
$ws.Range("C17").Value = [double]"0"
$ws.Range("D17").Value = [double]"99.9999"

$ws.Range("B18").Value = [double]"0.625"
$ws.Range("C18").Value = [double]"0"
$ws.Range("D18").Value = [double]"99.9999"

$ws.Range("B19").Value = [double]"100000"
$ws.Range("C19").Value = [double]"7.34609E-17"
$ws.Range("D19").Value = [double]"0.331988"

$ws.Range("B20").Value = [double]"10000000"
$ws.Range("C20").Value = [double]"6.28408E-20"
$ws.Range("D20").Value = [double]"14.2806"

# Apply the scientific number format (numFmtId 11) to every numeric cell in
# the three data blocks, matching style index 1 in the target styles.xml.
$ws.Range("B5:D8").NumberFormat = "0.00E+00"
$ws.Range("B11:D14").NumberFormat = "0.00E+00"
$ws.Range("B17:D20").NumberFormat = "0.00E+00"
